$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("TYPE") values for every row of the dictionary sheet, keyed by row number.
$typeByRow = @{
    1 = "TYPE"
    2 = "ID"
    3 = "location_data"
    4 = "location_data"
    5 = "location_data"
    6 = "location_data"
    7 = "location_data"
    8 = "location_data"
    9 = "safety_measures"
    10 = "safety_measures"
    11 = "traffic_data"
    12 = "traffic_data"
    13 = "traffic_data"
    14 = "traffic_data"
    15 = "traffic_data"
    16 = "traffic_data"
    17 = "traffic_data"
    18 = "traffic_data"
    19 = "traffic_data"
    20 = "traffic_data"
    21 = "traffic_data"
    22 = "traffic_data"
    23 = "traffic_data"
    24 = "traffic_data"
    25 = "traffic_data"
    26 = "traffic_data"
    27 = "traffic_data"
    28 = "geometry_data"
    29 = "geometry_data"
    30 = "geometry_data"
    31 = "geometry_data"
    32 = "safety_measures"
    33 = "safety_measures"
    34 = "safety_measures"
    35 = "safety_measures"
    36 = "safety_measures"
    37 = "safety_measures"
    38 = "safety_measures"
    39 = "safety_measures"
    40 = "safety_measures"
    41 = "traffic_data"
    42 = "traffic_data"
    43 = "traffic_data"
    44 = "traffic_data"
    45 = "traffic_data"
    46 = "traffic_data"
    47 = "traffic_data"
    48 = "traffic_data"
    49 = "geometry_data"
    50 = "geometry_data"
    51 = "safety_measures"
    52 = "geometry_data"
    53 = "geometry_data"
    54 = "safety_measures"
    55 = "safety_measures"
    56 = "location_data"
    57 = "location_data"
    58 = "traffic_data"
    59 = "traffic_data"
    60 = "location_data"
}

# Excel appends a brand-new entry to the shared-string table the first time a
# distinct string is written to a cell. The reference workbook's shared
# strings show these six new labels were introduced in this exact order:
# TYPE, ID, traffic_data, geometry_data, safety_measures, location_data.
# Seed one row per label (in that order) before filling in the rest, so the
# shared-string table is rebuilt identically.
$seedOrder = @("TYPE", "ID", "traffic_data", "geometry_data", "safety_measures", "location_data")

$written = @{}
foreach ($label in $seedOrder) {
    for ($row = 1; $row -le 60; $row++) {
        if (-not $written.ContainsKey($row)) {
            if ($typeByRow[$row] -eq $label) {
                $ws.Cells.Item($row, 3).Value = $typeByRow[$row]
                $written[$row] = $true
                break
            }
        }
    }
}

for ($row = 1; $row -le 60; $row++) {
    if (-not $written.ContainsKey($row)) {
        $ws.Cells.Item($row, 3).Value = $typeByRow[$row]
        $written[$row] = $true
    }
}

# Header cell (C1) matches the bold style already used by A1/B1
[void]($ws.Range("C1").Font.Bold = $true)

# Column C width, sized to fit its longest value ("safety_measures")
[void]($ws.Columns("C").AutoFit())

# Scroll the view back to the top-left and keep the existing selection
[void]($ws.Range("A1").Select())
[void]($ws.Range("B55").Select())
